# Add a new worksheet "ODI Batting Extra" at the end of the workbook,
# mirroring the header style of the existing "ODI Batting" sheet, and
# populate it with the per-match batting-extras data.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("ODI Batting")

# Insert the new sheet after the last existing sheet so it lands at the end.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "ODI Batting Extra"

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# --- Data rows --------------------------------------------------------------
# MATCH_CODE, NUM_4, NUM_6 and PERCENT_RUNS_OF_TOTAL hold numeric-looking
# text, so those columns are pre-formatted as Text before the values are
# written, which keeps them as strings (matching the source data) instead of
# being auto-converted to numbers/percentages by Excel.
$textCols = @(1, 3, 4, 5)
foreach ($col in $textCols) {
  $ws.Range($ws.Cells.Item(2, $col), $ws.Cells.Item(7, $col)).NumberFormat = "@"
}

$data = @(
  @("4402", $null,    $null, $null, $null,     "NO"),
  @("4406", $null,    $null, $null, $null,     "NO"),
  @("4410", 1,        "3",   "2",   "13.51%",  "NO"),
  @("4480", 1,        "9",   "0",   "16.35%",  "YES"),
  @("4482", 1,        "3",   "0",   "4.69%",   "NO"),
  @("4485", 1,        "8",   "0",   "21.78%",  "NO")
)

$r = 2
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  if ($row[1] -ne $null) {
    $ws.Cells.Item($r, 2).Value = $row[1]
  }
  if ($row[2] -ne $null) {
    $ws.Cells.Item($r, 3).Value = $row[2]
  }
  if ($row[3] -ne $null) {
    $ws.Cells.Item($r, 4).Value = $row[3]
  }
  if ($row[4] -ne $null) {
    $ws.Cells.Item($r, 5).Value = $row[4]
  }
  $ws.Cells.Item($r, 6).Value = $row[5]
  $r = $r + 1
}

# Drop the temporary "@" text format now that the literal strings are
# committed, so the data cells fall back to the default (unstyled) look,
# same as the source workbook's data rows.
$ws.Range("A2:F7").ClearFormats()

# Reuse the exact header formatting (bold font, border, centered) from the
# "ODI Batting" sheet's header row. Done last so ClearFormats above cannot
# touch it.
$src.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$null = $ws.Range("A1").Select()
